# Updated the GUI.excel_output to close when the crawling starts.
# Refresh the data table with a new crawl's results (new reference codes,
# a single new date, new credit amounts) and drop the last three rows of
# the previous crawl (the "Working Professional" subtotal rows), leaving
# just one "Working Professional" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop old row 14 (SM-3362BB) first: this is the last row of the
# "Student" merged block (A2:A14), so removing it shrinks the merge to
# A2:A13 without Excel re-splitting the merged block's borders, and it
# shifts the old "Working Professional" row (15) up into row 14 keeping
# its existing A-column label intact.
$ws.Range("A14:E14").EntireRow.Delete()

# The two remaining old rows (now at 15/16 after the shift above) are no
# longer needed - drop them too, leaving a single "Working Professional" row.
$ws.Range("A15:E16").EntireRow.Delete()

# Refresh the reference codes (column B) for the new crawl.
$ws.Range("B2").Value = "AS-CB728C"
$ws.Range("B3").Value = "CA-CB49D4"
$ws.Range("B4").Value = "CG-5EB883"
$ws.Range("B5").Value = "CG-CB74ED"
$ws.Range("B6").Value = "CL-BE26EF"
$ws.Range("B7").Value = "GH-FE0E73"
$ws.Range("B8").Value = "IM-74682B"
$ws.Range("B9").Value = "JE-CB64CF"
$ws.Range("B10").Value = "MS-CB7818"
$ws.Range("B11").Value = "RB-5EEEF0"
$ws.Range("B12").Value = "SR-3F1063"
$ws.Range("B13").Value = "YW-0DA6F3"
$ws.Range("B14").Value = "RA-044F6A"

# Every row now shares the same crawl date.
$ws.Range("C2:C14").Value = "22-07-2021"

# Refresh the credit amounts (column D).
$ws.Range("D2").Value = 391.04
$ws.Range("D3").Value = 690
$ws.Range("D4").Value = 443.44
$ws.Range("D5").Value = 461.14
$ws.Range("D6").Value = 311.76
$ws.Range("D7").Value = 295
$ws.Range("D8").Value = 404.8000000000001
$ws.Range("D9").Value = 398.4
$ws.Range("D10").Value = 436.16
$ws.Range("D11").Value = 478.69
$ws.Range("D12").Value = 347.76
$ws.Range("D13").Value = 658.16
$ws.Range("D14").Value = 977.5200000000002

# Refresh the Reference counts (column E).
$ws.Range("E2").Value = 8
$ws.Range("E3").Value = 8
$ws.Range("E4").Value = 8
$ws.Range("E5").Value = 1
$ws.Range("E6").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 8
$ws.Range("E9").Value = 8
$ws.Range("E10").Value = 8
$ws.Range("E11").Value = 1
$ws.Range("E12").Value = 8
$ws.Range("E13").Value = 1
$ws.Range("E14").Value = 8
